# Insert a new row at position 389 (pushes existing rows 389-408 down to 390-409)
# and populate it with the new weekly price observation.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(389).Insert()

$ws.Cells.Item(389, 1).Value2 = 4
$ws.Cells.Item(389, 2).Value2 = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(389, 3).Value2 = "Los Lagos"
$ws.Cells.Item(389, 4).Value2 = 45041
$ws.Cells.Item(389, 5).Value2 = 10
$ws.Cells.Item(389, 6).Value2 = 100112021
$ws.Cells.Item(389, 7).Value2 = "Ají"
$ws.Cells.Item(389, 8).Value2 = "Inferno"
$ws.Cells.Item(389, 9).Value2 = "Primera"
$ws.Cells.Item(389, 10).Value2 = 150
$ws.Cells.Item(389, 11).Value2 = 25000
$ws.Cells.Item(389, 12).Value2 = 25000
$ws.Cells.Item(389, 13).Value2 = 25000
$ws.Cells.Item(389, 14).Value2 = "`$/caja 10 kilos"
$ws.Cells.Item(389, 15).Value2 = "Región de Arica y Parinacota"
$ws.Cells.Item(389, 16).Value2 = 2500
$ws.Cells.Item(389, 17).Value2 = 10
$ws.Cells.Item(389, 18).Value2 = "Hortaliza"
